# fix(module3): use uncon_planned_qty for future production; keep produced for today
# Insert a new row (MAT_A / PLANT_001) after row 3, shifting the existing
# MAT_B rows down by one, and update the quantity / layer / horizon values
# for every data row to reflect the corrected calculation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 4:6 down to 5:7 to make room for the new MAT_A/PLANT_001 row,
# preserving the existing row/cell formatting.
$ws.Rows.Item(4).Insert()

# Row 2 - MAT_A / DC_001
$ws.Range("F2").Value = -239

# Row 3 - MAT_A / DC_002
$ws.Range("F3").Value = -562

# Row 4 (new) - MAT_A / PLANT_001
$ws.Range("A4").Value = "MAT_A"
$ws.Range("B4").Value = "PLANT_001"
$ws.Range("C4").Value = 45298
$ws.Range("C4").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D4").Value = "Distribution Demand - Forecast"
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = -863
$ws.Range("G4").Value = 45297
$ws.Range("G4").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("H4").Value = 1

# Row 5 - MAT_B / DC_001 (was row 4)
$ws.Range("B5").Value = "DC_001"
$ws.Range("F5").Value = -113
$ws.Range("H5").Value = 4

# Row 6 - MAT_B / DC_002 (was row 5)
$ws.Range("B6").Value = "DC_002"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = -32
$ws.Range("H6").Value = 1

# Row 7 - MAT_B / PLANT_001 (was row 6)
$ws.Range("A7").Value = "MAT_B"
$ws.Range("B7").Value = "PLANT_001"
$ws.Range("C7").Value = 45298
$ws.Range("C7").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D7").Value = "Distribution Demand - Forecast"
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = -113
$ws.Range("G7").Value = 45297
$ws.Range("G7").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("H7").Value = 1
